# Week 13 logging update
#  - Rushing sheet: new player "G.Minshew" inserted as row 3 (pushing the
#    rest of the roster down one row), plus updated cumulative rushing
#    stats for the week.
#  - Receiving sheet: updated cumulative receiving stats for the week
#    (no new rows there).

$wb = $excel.ActiveWorkbook
$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet: insert a new row for G.Minshew at row 3 -----------------
[void]$rushing.Rows("3:3").Insert(-4121)

# Bring over the formatting (border/alignment on col A) from the row above
# so the new row matches the rest of the table.
$rushing.Range("A2:F2").Copy()
[void]$rushing.Range("A3:F3").PasteSpecial(-4122)

# --- Rushing sheet: rewrite the full data block with this week's numbers ---
$rushingData = @(
    @(0, "J.Hurts",    41, 38, 33, 26),
    @(1, "G.Minshew",   2,  0,  2,  0),
    @(2, "M.Sanders",  30, 17,  2,  5),
    @(3, "B.Scott",    46, 24,  8, 14),
    @(4, "K.Gainwell", 43, 30,  7, 14),
    @(5, "J.Howard",   24, 20,  7, 17),
    @(6, "J.Reagor",    5,  1,  0,  0),
    @(7, "Q.Watkins",   1,  0,  0,  0),
    @(8, "M.Walker",    0,  0,  1,  0)
)

$r = 2
foreach ($row in $rushingData) {
    $rushing.Cells.Item($r, 1).Value = $row[0]
    $rushing.Cells.Item($r, 2).Value = $row[1]
    $rushing.Cells.Item($r, 3).Value = $row[2]
    $rushing.Cells.Item($r, 4).Value = $row[3]
    $rushing.Cells.Item($r, 5).Value = $row[4]
    $rushing.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Receiving sheet: rewrite the full data block with this week's numbers -
$receivingData = @(
    @(0, "M.Sanders",            4,  4,  1,  0, 0, 0),
    @(1, "B.Scott",              12,  9,  0,  0, 1, 0),
    @(2, "K.Gainwell",           39, 28,  2,  0, 8, 7),
    @(3, "D.Smith",              54, 39, 27, 11, 5, 3),
    @(4, "J.Reagor",             35, 23, 11,  3, 4, 2),
    @(5, "Q.Watkins",            28, 21, 17,  9, 8, 4),
    @(6, "G.Ward",                5,  2,  0,  0, 5, 2),
    @(7, "J.Arcega-Whiteside",    1,  1,  1,  1, 0, 0),
    @(8, "D.Goedert",            41, 29, 15, 12, 7, 5),
    @(9, "J.Stoll",                3,  2,  0,  0, 1, 1),
    @(10, "T.Jackson",             1,  0,  0,  0, 0, 0)
)

$r = 2
foreach ($row in $receivingData) {
    $receiving.Cells.Item($r, 1).Value = $row[0]
    $receiving.Cells.Item($r, 2).Value = $row[1]
    $receiving.Cells.Item($r, 3).Value = $row[2]
    $receiving.Cells.Item($r, 4).Value = $row[3]
    $receiving.Cells.Item($r, 5).Value = $row[4]
    $receiving.Cells.Item($r, 6).Value = $row[5]
    $receiving.Cells.Item($r, 7).Value = $row[6]
    $receiving.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

Write-Host "Week 13 logging complete"
